$wb = $excel.ActiveWorkbook

# The "remove_vars" sheet is being split: a new sheet "remove_vars_contact"
# is created holding the first 120 variable names (contact-related vars),
# placed right after "remove_vars" and becoming the active tab.
$values = @(
    'remove',
    'bestdregion ',
    'cultureinfo ',
    'ical_standardoccupation ',
    'insq61a ',
    'insq61b ',
    'main_intro_screen ',
    'nameagesex_scale_1 ',
    'nostdregion ',
    'part_ukitv ',
    'privacypolicy_insert ',
    'q_panel ',
    'q_wave ',
    'qconfirm1a ',
    'qconfirm1b ',
    'qedu_no ',
    'qforigin_codes ',
    'qforigin_rec ',
    'qmorigin_codes ',
    'qmorigin_rec ',
    'qn086 ',
    'qno81 ',
    'qrorigin_codes ',
    'qrorigin_rec ',
    'qsample ',
    'quotagerange ',
    'sniffer_device_type_final ',
    'sniffer_device_type_initial ',
    'uk02eth_18 ',
    'ukmktsize ',
    'ukmktsize2 ',
    'ukmunicipality ',
    'uksg_version ',
    'usegendernonbinary ',
    'w1_namesexage ',
    'w2_namesexage ',
    'w3_namesexage ',
    'w4_namesexage ',
    'w5_namesexage ',
    'w6_namesexage ',
    'qmktsize_refused_to_answer ',
    'qmktsize_text ',
    'qmktsizeversion ',
    'qmktsize_15_1 ',
    'qmktsize_16_1 ',
    'qmktsize_32_1 ',
    'qmktsize_33_1 ',
    'qmktsize_44_1 ',
    'qmktsize_45_1 ',
    'qmktsize_46_1 ',
    'qmktsize_47_1 ',
    'qmktsize_48_1 ',
    'qmktsize_49_1 ',
    'qmktsize_50_1 ',
    'qmktsize_51_1 ',
    'qmktsize_52_1 ',
    'qmktsize_53_1 ',
    'qmktsize_54_1 ',
    'qmktsize_55_1 ',
    'qmktsize_57_1 ',
    'qmktsize_58_1 ',
    'qmktsize_59_1 ',
    'qmktsize_60_1 ',
    'qmktsize_61_1 ',
    'qmktsize_62_1 ',
    'qmktsize_63_1 ',
    'qmktsize_64_1 ',
    'qmktsize_65_1 ',
    'qmktsize_66_1 ',
    'qmktsize_67_1 ',
    'qmktsize_68_1 ',
    'qmktsize_69_1 ',
    'qmktsize_70_1 ',
    'qmktsize_71_1 ',
    'qmktsize_72_1 ',
    'qmktsize_73_1 ',
    'qmktsize_74_1 ',
    'qmktsize_75_1 ',
    'qmktsize_76_1 ',
    'qmktsize_77_1 ',
    'qmktsize_78_1 ',
    'qmktsize_79_1 ',
    'qmktsize_80_1 ',
    'qmktsize_81_1 ',
    'qmktsize_82_1 ',
    'qmktsize_19_1 ',
    'qmktsize_85_1 ',
    'qmktsize_86_1 ',
    'qmktsize_34_1 ',
    'qmktsize_35_1 ',
    'qmktsize_36_1 ',
    'qmktsize_37_1',
    'q20 ',
    'q20_import ',
    'q20_new ',
    'q20_original ',
    'q20b ',
    'q36_insert_4 ',
    'q39 ',
    'q40 ',
    'q41 ',
    'q48a ',
    'q49_l0 ',
    'q49_l1 ',
    'q53a_choosen ',
    'q61a_insert ',
    'q61b_codes_1 ',
    'q61b_codes_2 ',
    'q61b_codes_98 ',
    'q61b_codes_99 ',
    'q61b_insert ',
    'q63validate ',
    'uk_nations ',
    'uk_region3_label ',
    'uk02inc ',
    'ukregion1 ',
    'ukruralurban ',
    'ukstdregion ',
    'w7_namesexage ',
    'w8_namesexage '
)

$sourceSheet = $wb.Worksheets.Item("remove_vars")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sourceSheet)
$newSheet.Name = "remove_vars_contact"

for ($i = 0; $i -lt $values.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $values[$i]
}

$newSheet.Activate()
